$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Clean up posting_match (M) title lists: move mis-ordered titles so
#     they sit correctly relative to the entries they were parsed from ---
$ws.Cells.Item(4, 13).Value = "比部員外郎;給事中;宣州;翰林學士;判太常卿事;水部員外郎;衛尉卿;營田副使;江州;右諫議大夫;御史中丞;知制誥;司士參軍;江州;主客郎中"
$ws.Cells.Item(8, 13).Value = "祠部郎中;大理評事;殿中監;給事中;駕部員外郎;秘書省正字;水部員外郎;司農少卿;屯田員外郎;修國史;刑部侍郎;知制誥;中書舍人"

# --- Row 22: a writing (T/U) match was found after cleaning the title text ---
$ws.Cells.Item(22, 4).Value = 12
$ws.Cells.Item(22, 20).Value = 1
$ws.Cells.Item(22, 21).Value = "文集"

# --- Force column A (input_id) to stay text for the new rows, matching the
#     existing sheet where numeric-looking ids are stored as strings ---
$idRange = $ws.Range("A24:A26")
$idRange.NumberFormat = "@"

# --- New rows found via the cleaned writing-title match ---
$ws.Cells.Item(24, 1).Value = "11833"
$ws.Cells.Item(24, 2).Value = 11833
$ws.Cells.Item(24, 3).Value = "劉玶"
$ws.Cells.Item(24, 4).Value = 1
$ws.Cells.Item(24, 5).Value = 0
$ws.Cells.Item(24, 7).Value = 0
$ws.Cells.Item(24, 9).Value = 0
$ws.Cells.Item(24, 11).Value = "有詩集存世"
$ws.Cells.Item(24, 12).Value = 0
$ws.Cells.Item(24, 14).Value = 0
$ws.Cells.Item(24, 16).Value = 0
$ws.Cells.Item(24, 18).Value = 0
$ws.Cells.Item(24, 20).Value = 1
$ws.Cells.Item(24, 21).Value = "詩集"
$ws.Cells.Item(24, 22).Value = "15;宋"
$ws.Cells.Item(24, 23).Value = "no"

$ws.Cells.Item(25, 1).Value = "10111"
$ws.Cells.Item(25, 2).Value = 10111
$ws.Cells.Item(25, 3).Value = "唐慎微"
$ws.Cells.Item(25, 4).Value = 1
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 7).Value = 0
$ws.Cells.Item(25, 9).Value = 0
$ws.Cells.Item(25, 11).Value = "重修政和經史證類備用本草"
$ws.Cells.Item(25, 12).Value = 0
$ws.Cells.Item(25, 14).Value = 0
$ws.Cells.Item(25, 16).Value = 0
$ws.Cells.Item(25, 18).Value = 0
$ws.Cells.Item(25, 20).Value = 1
$ws.Cells.Item(25, 21).Value = "經史證類備用本草"
$ws.Cells.Item(25, 22).Value = "15;宋"
$ws.Cells.Item(25, 23).Value = "no"

$ws.Cells.Item(26, 1).Value = "10831"
$ws.Cells.Item(26, 2).Value = 10831
$ws.Cells.Item(26, 3).Value = "李心傳"
$ws.Cells.Item(26, 4).Value = 1
$ws.Cells.Item(26, 5).Value = 0
$ws.Cells.Item(26, 7).Value = 0
$ws.Cells.Item(26, 9).Value = 0
$ws.Cells.Item(26, 11).Value = "著建炎以來朝野雜記"
$ws.Cells.Item(26, 12).Value = 0
$ws.Cells.Item(26, 14).Value = 0
$ws.Cells.Item(26, 16).Value = 0
$ws.Cells.Item(26, 18).Value = 0
$ws.Cells.Item(26, 20).Value = 1
$ws.Cells.Item(26, 21).Value = "以來朝野雜記"
$ws.Cells.Item(26, 22).Value = "15;宋"
$ws.Cells.Item(26, 23).Value = "no"

# --- Mark the remaining empty/text-typed columns on the new rows as empty
#     inline strings (rather than leaving the cells absent), matching the
#     sheet's convention for "no match" columns ---
$emptyRange = $ws.Range("F24:F26,H24:H26,J24:J26,M24:M26,O24:O26,Q24:Q26,S24:S26")
$emptyRange.NumberFormat = "@"
foreach ($r in 24..26) {
    $ws.Cells.Item($r, 6).Value = " "
    $ws.Cells.Item($r, 8).Value = " "
    $ws.Cells.Item($r, 10).Value = " "
    $ws.Cells.Item($r, 13).Value = " "
    $ws.Cells.Item($r, 15).Value = " "
    $ws.Cells.Item($r, 17).Value = " "
    $ws.Cells.Item($r, 19).Value = " "
}
foreach ($r in 24..26) {
    $ws.Cells.Item($r, 6).Value = ""
    $ws.Cells.Item($r, 8).Value = ""
    $ws.Cells.Item($r, 10).Value = ""
    $ws.Cells.Item($r, 13).Value = ""
    $ws.Cells.Item($r, 15).Value = ""
    $ws.Cells.Item($r, 17).Value = ""
    $ws.Cells.Item($r, 19).Value = ""
}

$idRange.Style = "Normal"
$emptyRange.Style = "Normal"
